$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2612.111
$ws.Cells.Item(32, 10).Value = 2278.3333
$ws.Cells.Item(32, 12).Value = 2278.3333
$ws.Cells.Item(32, 14).Value = -2930.3333

$ws.Cells.Item(51, 8).Value = 4650.8
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 4650.8
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 4650.8
$ws.Cells.Item(51, 13).ClearContents()
$ws.Cells.Item(51, 14).Value = -5618.8

$ws.Cells.Item(96, 8).Value = 408.66666
$ws.Cells.Item(96, 9).Value = 397.72726
$ws.Cells.Item(96, 10).Value = 529
$ws.Cells.Item(96, 11).Value = 1193.18178
$ws.Cells.Item(96, 12).Value = 1587
$ws.Cells.Item(96, 13).Value = 179.8182200000001
$ws.Cells.Item(96, 14).Value = -4333

$ws.Cells.Item(97, 8).Value = 4996.6665
$ws.Cells.Item(97, 10).Value = 4996.6665
$ws.Cells.Item(97, 12).Value = 14989.9995
$ws.Cells.Item(97, 14).Value = -15981.9995

$ws.Cells.Item(112, 8).Value = 1910.2667
$ws.Cells.Item(112, 10).Value = 2042.9231
$ws.Cells.Item(112, 12).Value = 6128.7693
$ws.Cells.Item(112, 14).Value = -8344.7693

$ws.Cells.Item(136, 8).Value = 77973.75
$ws.Cells.Item(136, 10).Value = 77973.75
$ws.Cells.Item(136, 12).Value = 77973.75
$ws.Cells.Item(136, 14).Value = -88173.75

$ws.Cells.Item(137, 8).Value = 2450.9148
$ws.Cells.Item(137, 9).Value = 1528.15
$ws.Cells.Item(137, 10).Value = 3134.4443
$ws.Cells.Item(137, 11).Value = 4584.450000000001
$ws.Cells.Item(137, 12).Value = 9403.332900000001
$ws.Cells.Item(137, 13).Value = -2034.450000000001
$ws.Cells.Item(137, 14).Value = -14503.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2276728.8
$ws.Cells.Item(74, 9).Value = 2844255.5
$ws.Cells.Item(74, 10).Value = 6620.909
$ws.Cells.Item(74, 11).Value = 2844255.5
$ws.Cells.Item(74, 12).Value = 6620.909
$ws.Cells.Item(74, 13).Value = -2843381.5
$ws.Cells.Item(74, 14).Value = -8368.909

$ws.Cells.Item(77, 8).Value = 2276728.8
$ws.Cells.Item(77, 9).Value = 2844255.5
$ws.Cells.Item(77, 10).Value = 6620.909
$ws.Cells.Item(77, 11).Value = 14221277.5
$ws.Cells.Item(77, 12).Value = 33104.545
$ws.Cells.Item(77, 13).Value = -14216909.5
$ws.Cells.Item(77, 14).Value = -41840.545

$ws.Cells.Item(117, 8).Value = 80248
$ws.Cells.Item(117, 10).Value = 80248
$ws.Cells.Item(117, 12).Value = 80248
$ws.Cells.Item(117, 14).Value = -89426

$ws.Cells.Item(132, 8).Value = 451668.12
$ws.Cells.Item(132, 9).Value = 544243.0600000001
$ws.Cells.Item(132, 11).Value = 1632729.18
$ws.Cells.Item(132, 13).Value = -1630199.18

$ws.Cells.Item(133, 8).Value = 69999
$ws.Cells.Item(133, 10).Value = 69999
$ws.Cells.Item(133, 12).Value = 69999
$ws.Cells.Item(133, 14).Value = -75059

$ws.Cells.Item(134, 8).Value = 85715.5
$ws.Cells.Item(134, 10).Value = 85715.5
$ws.Cells.Item(134, 12).Value = 85715.5
$ws.Cells.Item(134, 14).Value = -95855.5

$ws.Cells.Item(137, 8).Value = 88000
$ws.Cells.Item(137, 10).Value = 88000
$ws.Cells.Item(137, 12).Value = 88000
$ws.Cells.Item(137, 14).Value = -98200

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 79780
$ws.Cells.Item(52, 10).Value = 79780
$ws.Cells.Item(52, 12).Value = 79780
$ws.Cells.Item(52, 14).Value = -80306

$ws.Cells.Item(121, 8).Value = 79780
$ws.Cells.Item(121, 10).Value = 79780
$ws.Cells.Item(121, 12).Value = 79780
$ws.Cells.Item(121, 14).Value = -83274

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 12254.546
$ws.Cells.Item(35, 9).Value = 13542.857
$ws.Cells.Item(35, 11).Value = 13542.857
$ws.Cells.Item(35, 13).Value = -13248.857

$ws.Cells.Item(70, 8).Value = 75262
$ws.Cells.Item(70, 9).Value = 74999
$ws.Cells.Item(70, 10).Value = 75393.5
$ws.Cells.Item(70, 11).Value = 74999
$ws.Cells.Item(70, 12).Value = 75393.5
$ws.Cells.Item(70, 13).Value = -74684
$ws.Cells.Item(70, 14).Value = -76023.5

$ws.Cells.Item(73, 8).Value = 75262
$ws.Cells.Item(73, 9).Value = 74999
$ws.Cells.Item(73, 10).Value = 75393.5
$ws.Cells.Item(73, 11).Value = 74999
$ws.Cells.Item(73, 12).Value = 75393.5
$ws.Cells.Item(73, 13).Value = -73907
$ws.Cells.Item(73, 14).Value = -77577.5

$ws.Cells.Item(134, 8).Value = 11036.837
$ws.Cells.Item(134, 9).Value = 9460
$ws.Cells.Item(134, 11).Value = 28380
$ws.Cells.Item(134, 13).Value = -25845

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 63468664
$ws.Cells.Item(4, 9).Value = 102890120
$ws.Cells.Item(4, 11).Value = 308670360
$ws.Cells.Item(4, 13).Value = -308670248

$ws.Cells.Item(5, 8).Value = 2754.276
$ws.Cells.Item(5, 10).Value = 5177.5713
$ws.Cells.Item(5, 12).Value = 15532.7139
$ws.Cells.Item(5, 14).Value = -15756.7139

$ws.Cells.Item(16, 8).Value = 1497.1666
$ws.Cells.Item(16, 10).Value = 2797
$ws.Cells.Item(16, 12).Value = 8391
$ws.Cells.Item(16, 14).Value = -8737

$ws.Cells.Item(17, 8).Value = 174.4762
$ws.Cells.Item(17, 10).Value = 996.3333
$ws.Cells.Item(17, 12).Value = 2988.9999
$ws.Cells.Item(17, 14).Value = -3326.9999

$ws.Cells.Item(18, 8).Value = 1020
$ws.Cells.Item(18, 9).Value = 794.8570999999999
$ws.Cells.Item(18, 10).Value = 1414
$ws.Cells.Item(18, 11).Value = 2384.5713
$ws.Cells.Item(18, 12).Value = 4242
$ws.Cells.Item(18, 13).Value = -2215.5713
$ws.Cells.Item(18, 14).Value = -4580

$ws.Cells.Item(21, 8).Value = 3996.3333
$ws.Cells.Item(21, 9).Value = 1994.5
$ws.Cells.Item(21, 10).Value = 8000
$ws.Cells.Item(21, 11).Value = 5983.5
$ws.Cells.Item(21, 12).Value = 24000
$ws.Cells.Item(21, 13).Value = -5810.5
$ws.Cells.Item(21, 14).Value = -24346

$ws.Cells.Item(54, 8).Value = 5666.6665
$ws.Cells.Item(54, 10).Value = 5666.6665
$ws.Cells.Item(54, 12).Value = 16999.9995
$ws.Cells.Item(54, 14).Value = -18117.9995

$ws.Cells.Item(76, 8).Value = 6250.75
$ws.Cells.Item(76, 9).Value = 4506.5
$ws.Cells.Item(76, 11).Value = 13519.5
$ws.Cells.Item(76, 13).Value = -13136.5

$ws.Cells.Item(79, 8).Value = 6250.75
$ws.Cells.Item(79, 9).Value = 4506.5
$ws.Cells.Item(79, 11).Value = 13519.5
$ws.Cells.Item(79, 13).Value = -12193.5

$ws.Cells.Item(93, 8).Value = 5810.3
$ws.Cells.Item(93, 10).Value = 6011.4443
$ws.Cells.Item(93, 12).Value = 18034.3329
$ws.Cells.Item(93, 14).Value = -21778.3329

$ws.Cells.Item(99, 8).Value = 8922.4
$ws.Cells.Item(99, 10).Value = 8922.4
$ws.Cells.Item(99, 12).Value = 26767.2
$ws.Cells.Item(99, 14).Value = -31259.2

$ws.Cells.Item(100, 8).Value = 2666
$ws.Cells.Item(100, 10).Value = 2499.5
$ws.Cells.Item(100, 12).Value = 7498.5
$ws.Cells.Item(100, 14).Value = -9120.5

$ws.Cells.Item(101, 8).Value = 7903.75
$ws.Cells.Item(101, 10).Value = 7903.75
$ws.Cells.Item(101, 12).Value = 23711.25
$ws.Cells.Item(101, 14).Value = -28579.25

$ws.Cells.Item(107, 8).Value = 820.1818
$ws.Cells.Item(107, 10).Value = 917
$ws.Cells.Item(107, 12).Value = 2751
$ws.Cells.Item(107, 14).Value = -6591

$ws.Cells.Item(109, 8).Value = 4782.647
$ws.Cells.Item(109, 10).Value = 5904.5454
$ws.Cells.Item(109, 12).Value = 17713.6362
$ws.Cells.Item(109, 14).Value = -19793.6362

$ws.Cells.Item(115, 8).Value = 9000
$ws.Cells.Item(115, 9).Value = 9000
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 27000
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = -25825
$ws.Cells.Item(115, 14).ClearContents()

$ws.Cells.Item(117, 8).Value = 3836.1304
$ws.Cells.Item(117, 9).Value = 1724.4286
$ws.Cells.Item(117, 10).Value = 4760
$ws.Cells.Item(117, 11).Value = 5173.2858
$ws.Cells.Item(117, 12).Value = 14280
$ws.Cells.Item(117, 13).Value = -1731.2858
$ws.Cells.Item(117, 14).Value = -21164

$ws.Cells.Item(118, 8).Value = 1587.8334
$ws.Cells.Item(118, 10).Value = 1400
$ws.Cells.Item(118, 12).Value = 4200
$ws.Cells.Item(118, 14).Value = -6686

$ws.Cells.Item(119, 8).Value = 700
$ws.Cells.Item(119, 9).Value = 700
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = 2100
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 13).Value = 2738
$ws.Cells.Item(119, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 539.7742
$ws.Cells.Item(122, 9).Value = 483
$ws.Cells.Item(122, 10).Value = 550.6923
$ws.Cells.Item(122, 11).Value = 4347
$ws.Cells.Item(122, 12).Value = 4956.2307
$ws.Cells.Item(122, 13).Value = -1897
$ws.Cells.Item(122, 14).Value = -9856.2307

$ws.Cells.Item(135, 8).Value = 2754.276
$ws.Cells.Item(135, 10).Value = 5177.5713
$ws.Cells.Item(135, 12).Value = 46598.14169999999
$ws.Cells.Item(135, 14).Value = -51668.14169999999

$ws.Cells.Item(136, 8).Value = 3088.1667
$ws.Cells.Item(136, 9).Value = 3088.1667
$ws.Cells.Item(136, 11).Value = 9264.500100000001
$ws.Cells.Item(136, 13).Value = -4164.500100000001

$ws.Cells.Item(138, 8).Value = 2261.8572
$ws.Cells.Item(138, 9).Value = 1972.1666
$ws.Cells.Item(138, 11).Value = 5916.4998
$ws.Cells.Item(138, 13).Value = -776.4997999999996

$ws.Cells.Item(139, 8).Value = 950
$ws.Cells.Item(139, 9).Value = 950
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 2850
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = 2290
$ws.Cells.Item(139, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(131, 8).Value = 27903.715
$ws.Cells.Item(131, 10).Value = 27903.715
$ws.Cells.Item(131, 12).Value = 27903.715
$ws.Cells.Item(131, 14).Value = -37983.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3728.6333
$ws.Cells.Item(122, 9).Value = 3531.5652
$ws.Cells.Item(122, 11).Value = 10594.6956
$ws.Cells.Item(122, 13).Value = -8144.695599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2702.3044
$ws.Cells.Item(122, 9).Value = 2459.6667
$ws.Cells.Item(122, 11).Value = 7379.000100000001
$ws.Cells.Item(122, 13).Value = -4929.000100000001
